# STEEXP1.xlsx - "SEM 5 - 28/11/2023"
#
# The sheet's test-case table (columns B:I) was re-laid-out: the data
# columns are auto-fitted to their (now longer) wrapped text, the
# "INPUT DATA" column (F) is given a bit of extra manual breathing room on
# top of that, and the active selection is left on F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Auto-fit every data column of the table to its contents.
$ws.Columns("B:I").AutoFit()

# Column F ("INPUT DATA") was then nudged a little wider by hand, so it no
# longer sits at its pure auto-fit width.
$ws.Columns("F:F").ColumnWidth = 24.83

# Leave the selection where the editor left it.
$ws.Range("F5").Select()
